$d = $word.ActiveDocument

$pairs = @(
    @("82÷7=", "68÷7="),
    @("51÷2=", "19÷8="),
    @("16÷4=", "67÷5="),
    @("49÷9=", "71÷7="),
    @("86÷8=", "14÷7="),
    @("32÷7=", "64÷3="),
    @("25÷4=", "30÷9="),
    @("63÷7=", "26÷6="),
    @("74÷6=", "94÷8="),
    @("61÷2=", "82÷5="),
    @("61÷4=", "36÷7="),
    @("98÷8=", "53÷2="),
    @("56÷5=", "47÷2="),
    @("48÷8=", "39÷9="),
    @("33÷7=", "10÷3="),
    @("48÷5=", "31÷2="),
    @("66÷2=", "93÷6="),
    @("86÷6=", "28÷9="),
    @("31÷3=", "65÷9="),
    @("92÷5=", "13÷4="),
    @("67÷8=", "50÷3="),
    @("38÷7=", "39÷9="),
    @("79÷8=", "54÷4="),
    @("74÷4=", "17÷3="),
    @("99÷8=", "34÷3=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
